$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Digital Print F 4x0 / 252-HP 10000 Press / 719" activity row (originally
# row 2) moves down to row 4, the "Cut / 406-45" Polar 115ED Cutter / 715" row
# (originally row 3) moves up to row 2, and the "- / 169-Press Approval Task /
# 740" row (originally row 4) moves up to row 3. Text cells are reassigned with
# .Value (kept as shared-string text); the purely-numeric PlannedQty cells (G)
# use .Replace so Excel does not silently reinterpret the text as a number and
# re-key the cell's style.

# Row 2 <- old row 3 (Cut)
$ws.Range("B2").Value = "Cut"
$ws.Range("D2").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("G2").Replace("719", "715")
$ws.Range("L2").Value = "406-45`" Polar 115ED Cutter"
$ws.Range("M2").Value = "406-45`" Polar 115ED Cutter`n404-45`" Polar 115EMC Cutter`n405-54`" Polar 137EMC Cutter`n402-45`" Polar 115EMC Cutter`n403-54`" Polar 137ED Cutter"

# Row 3 <- old row 4 (Press Approval Task)
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "169-Press Approval Task "
$ws.Range("G3").Replace("715", "740")
$ws.Range("L3").Value = "Press Approval Task"
$ws.Range("M3").Value = "169-Press Approval Task "

# Row 4 <- old row 2 (Digital Print F 4x0)
$ws.Range("B4").Value = "Digital Print F 4x0"
$ws.Range("D4").Value = "252-HP 10000 Press"
$ws.Range("G4").Replace("740", "719")
$ws.Range("L4").Value = "252-HP 10000 Press"
$ws.Range("M4").Value = "252-HP 10000 Press"
